# Applies the data refresh captured in the commit:
#  - fixes mojibake (mis-encoded accented characters) in the Regional
#    Economic Communities footnote
#  - updates a batch of recalculated aggregate statistics in rows 61-98
#    of the Tab12 worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footnote text correction (cell A103) ---------------------------------
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# --- Recalculated numeric cells -------------------------------------------
$ws.Range("C61").Value = 67.358333333333405
$ws.Range("J62").Value = 60.569696969696999
$ws.Range("D63").Value = 60.753571428571398
$ws.Range("G64").Value = 98.73
$ws.Range("C66").Value = 74.1142857142857
$ws.Range("F66").Value = 40.028571428571396
$ws.Range("G67").Value = 90.314285714285703
$ws.Range("C77").Value = 91.383333333333297
$ws.Range("H82").Value = 76.453333333333404
$ws.Range("J82").Value = 61.3466666666667
$ws.Range("H95").Value = 75.033333333333402
$ws.Range("J95").Value = 59.033333333333402
$ws.Range("G96").Value = 92.618181818181895

$ws.Range("C97").Value = 66.330434782608705
$ws.Range("D97").Value = 42.991304347826102
$ws.Range("E97").Value = 20.4434782608696
$ws.Range("F97").Value = 29.539130434782599
$ws.Range("G97").Value = 88.714285714285694
$ws.Range("H97").Value = 73.109523809523793
$ws.Range("I97").Value = 49.5
$ws.Range("J97").Value = 58.257142857142902
$ws.Range("K97").Value = 2093.8479029999999
$ws.Range("L97").Value = 10.409098627500001
$ws.Range("M97").Value = 5752.4042406999997
$ws.Range("N97").Value = 26.224801780958298
$ws.Range("O97").Value = 26.1884615384615

$ws.Range("C98").Value = 79.616666666666703
$ws.Range("D98").Value = 52.683333333333302
$ws.Range("E98").Value = 27.05
$ws.Range("F98").Value = 40.5
$ws.Range("G98").Value = 93.8333333333333
$ws.Range("H98").Value = 72.816666666666706
$ws.Range("I98").Value = 53.016666666666701
$ws.Range("J98").Value = 65.016666666666694
$ws.Range("K98").Value = 3202.1190000000001
$ws.Range("L98").Value = 7.608918708
$ws.Range("M98").Value = 8769.7732445000001
$ws.Range("N98").Value = 22.5908406223333
$ws.Range("O98").Value = 38.590909090909101
